# Add a new "Swiss" worksheet (copy of "Czech") with Switzerland market test data.

$wb = $excel.ActiveWorkbook

# Czech sheet is the template for the new Swiss sheet.
$czech = $wb.Worksheets.Item("Czech")

# Select the whole sheet on Czech (leftover selection state from copying it),
# matching the post-copy selection state of the source sheet.
$czech.Cells.Select()

# Duplicate "Czech" and move the copy to the end of the workbook, then rename it.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$czech.Copy($null, $lastSheet)
$swiss = $wb.Worksheets.Item($wb.Worksheets.Count)
$swiss.Name = "Swiss"

# Update the market-specific cells for Switzerland.
$swiss.Range("B2").Value = "Switzerland Market"
$swiss.Range("B4").Value = "NGC-3476/T2646"

# Mirror the saved selection/active-cell state on the new sheet and make it active.
$swiss.Activate()
$swiss.Range("B2:B4").Select()
